$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.518.65"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.659.98"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'307.66"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'0.9983"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.3616"
$ws.Range("E7").Value = "  -2.98%  "
$ws.Range("D8").Value = "'47.65"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "'0.3251"
$ws.Range("E9").Value = "  -5.36%  "
$ws.Range("D10").Value = "'1.124"
$ws.Range("E10").Value = "  -4.85%  "
$ws.Range("D11").Value = "'0.06955"
$ws.Range("E11").Value = "  -6.45%  "
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "'5.894"
$ws.Range("E13").Value = "  -5.21%  "
$ws.Range("D14").Value = "'19.38"
$ws.Range("E14").Value = "  -6.96%  "
$ws.Range("D15").Value = "1.655.16"
$ws.Range("E15").Value = "  -3.03%  "
$ws.Range("D16").Value = "'6.554"
$ws.Range("E16").Value = "  -5.26%  "
$ws.Range("D17").Value = "'0.00001045"
$ws.Range("E17").Value = "  -6.45%  "
$ws.Range("D18").Value = "'0.06543"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "'0.9979"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'76.46"
$ws.Range("E20").Value = "  -8.11%  "
$ws.Range("D21").Value = "'5.913"
$ws.Range("E21").Value = "  -6.63%  "
$ws.Range("D22").Value = "'15.65"
$ws.Range("E22").Value = "  -8.29%  "
$ws.Range("D23").Value = "'12.55"
$ws.Range("E23").Value = "  -4.55%  "
$ws.Range("D24").Value = "24.532.14"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").Value = "'2.469"
$ws.Range("E25").Value = "  +2.40%  "
$ws.Range("D26").Value = "'2.302"
$ws.Range("E26").Value = "  -16.43%  "
$ws.Range("D27").Value = "'146.79"
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").Value = "'18.48"
$ws.Range("E28").Value = "  -7.89%  "
$ws.Range("D29").Value = "1.847.04"
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("D30").Value = "'1.200"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("D32").Value = "'4.055"
$ws.Range("D33").Value = "'5.625"
$ws.Range("E33").Value = "  -16.22%  "
$ws.Range("D34").Value = "'1.703"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("D35").Value = "'0.08363"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("D36").Value = "'12.34"
$ws.Range("E36").Value = "  -9.58%  "
$ws.Range("D37").Value = "'5.187"
$ws.Range("E37").Value = "  -5.57%  "
$ws.Range("D38").Value = "'0.06047"
$ws.Range("E38").Value = "  -7.14%  "
$ws.Range("D39").Value = "'0.2058"
$ws.Range("E39").Value = "  -6.91%  "
$ws.Range("D40").Value = "'1.204"
$ws.Range("E40").Value = "  -5.31%  "
$ws.Range("D41").Value = "'8.195"
$ws.Range("E41").Value = "  -7.94%  "
$ws.Range("D42").Value = "'0.02182"
$ws.Range("E42").Value = "  -7.67%  "
$ws.Range("D43").Value = "'0.9991"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").Value = "'0.5883"
$ws.Range("E44").Value = "  -7.79%  "
$ws.Range("D45").Value = "'3.735"
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").Value = "'12.73"
$ws.Range("E46").Value = "  -7.45%  "
$ws.Range("D47").Value = "'0.5586"
$ws.Range("E47").Value = "  -7.76%  "
$ws.Range("D48").Value = "'122.26"
$ws.Range("E48").Value = "  -5.21%  "
$ws.Range("D49").Value = "'1.937"
$ws.Range("E49").Value = "  -8.10%  "
$ws.Range("D50").Value = "'0.06921"
$ws.Range("E50").Value = "  -4.59%  "
$ws.Range("D51").Value = "'74.08"
$ws.Range("E51").Value = "  -6.05%  "
